# Updates cryptos list values (prices in column D, 1h volume % in column E,
# plus two coin-name/link/price/volume row swaps) per the Jun 3 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value while forcing text storage (matches the source file,
# which keeps numeric-looking price strings, e.g. "1.00", as text so that
# trailing zeros/leading zeros survive) and then restores the cell style so no
# stray number-format/style is left behind on the cell.
function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "69.146.28"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "3.810.08"
$ws.Range("E3").Value = "  +0.34%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.33%  "
Set-TextValue "D5" "633.64"
$ws.Range("E5").Value = "  +5.50%  "
Set-TextValue "D6" "165.25"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "3.807.83"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  +0.31%  "
Set-TextValue "D12" "6.63"
$ws.Range("E12").Value = "  +3.07%  "
Set-TextValue "D13" "0.0000250"
$ws.Range("E13").Value = "  +0.03%  "
Set-TextValue "D14" "35.93"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "4.449.25"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "3.778.11"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "69.111.69"
$ws.Range("E17").Value = "  +1.48%  "
Set-TextValue "D18" "18.13"
$ws.Range("E18").Value = "  -1.07%  "
Set-TextValue "D19" "7.13"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  -0.28%  "
Set-TextValue "D21" "466.52"
$ws.Range("E21").Value = "  +0.78%  "
Set-TextValue "D22" "9.64"
$ws.Range("E23").Value = "  +1.55%  "
Set-TextValue "D24" "0.0000152"
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D26" "2.16"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D27" "11.97"
$ws.Range("E27").Value = "  -0.92%  "
Set-TextValue "D28" "10.06"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "3.960.22"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  +0.32%  "
Set-TextValue "D33" "7.28"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("E34").Value = "  -0.44%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "3.755.13"
$ws.Range("E36").Value = "  +0.22%  "
Set-TextValue "D37" "9.02"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("E39").Value = "  +7.12%  "
$ws.Range("E40").Value = "  +5.15%  "
Set-TextValue "D41" "5.90"
$ws.Range("E41").Value = "  +1.90%  "
Set-TextValue "D42" "0.976"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +0.00%  "
Set-TextValue "D45" "157.16"
$ws.Range("E45").Value = "  +3.92%  "
Set-TextValue "D46" "0.300"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D47" "43.40"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D48" "1.42"
$ws.Range("E48").Value = "  +4.91%  "
Set-TextValue "D49" "46.84"
$ws.Range("E49").Value = "  -1.37%  "
Set-TextValue "D50" "1.91"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("E51").Value = "  +0.93%  "
